$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Range("B2").Value = -3.00978071494291
$ws1.Range("C2").Value = 0.294620451678279
$ws1.Range("B3").Value = 0.138978900167005
$ws1.Range("C3").Value = 0.170187449721332

$ws2 = $wb.Worksheets.Item("lognormal")
$ws2.Range("B2").Value = 2.48739055832804
$ws2.Range("C2").Value = 0.423446230224185
$ws2.Range("B3").Value = -1.07834058869686
$ws2.Range("C3").Value = 0.159202180743191

$ws3 = $wb.Worksheets.Item("llogis")
$ws3.Range("B2").Value = -2.3202820840382
$ws3.Range("C2").Value = 0.12396398295321
$ws3.Range("B3").Value = 0.686335169404983
$ws3.Range("C3").Value = 0.16539975502107

$ws4 = $wb.Worksheets.Item("gompertz")
$ws4.Range("B2").Value = -2.65972894717063
$ws4.Range("C2").Value = 0.189263107546043
$ws4.Range("B3").Value = -0.00246882176944569
$ws4.Range("C3").Value = 0.0180428086912718

$ws6 = $wb.Worksheets.Item("weibull cov")
$ws6.Range("A2").Value = 0.086801210547113
$ws6.Range("B2").Value = -0.0421769037283446
$ws6.Range("A3").Value = -0.0421769037283446
$ws6.Range("B3").Value = 0.0289637680426511

$ws7 = $wb.Worksheets.Item("lognormal cov")
$ws7.Range("A2").Value = 0.179306709891073
$ws7.Range("B2").Value = -0.0627180539399422
$ws7.Range("A3").Value = -0.0627180539399422
$ws7.Range("B3").Value = 0.0253453343533878

$ws8 = $wb.Worksheets.Item("llogis cov")
$ws8.Range("A2").Value = 0.0153670690696236
$ws8.Range("B2").Value = -0.00232816030093871
$ws8.Range("A3").Value = -0.00232816030093871
$ws8.Range("B3").Value = 0.02735707896103

$ws9 = $wb.Worksheets.Item("gompertz cov")
$ws9.Range("A2").Value = 0.0358205238779849
$ws9.Range("B2").Value = -0.00230345389205733
$ws9.Range("A3").Value = -0.00230345389205733
$ws9.Range("B3").Value = 0.000325542945469833
